$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 13065
$ws1.Range("F5").Value = 89
$ws1.Range("F6").Value = 100
$ws1.Range("F10").Value = 13032
$ws1.Range("F11").Value = 296
$ws1.Range("F12").Value = 549
$ws1.Range("F13").Value = 8733
$ws1.Range("F14").Value = 7763
$ws1.Range("F15").Value = 210
$ws1.Range("F18").Value = 133
$ws1.Range("F19").Value = 992
$ws1.Range("F20").Value = 10
$ws1.Range("F23").Value = 188
$ws1.Range("F24").Value = 336

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 13065
$ws4.Range("F6").Value = 89
$ws4.Range("F7").Value = 100
$ws4.Range("F11").Value = 13032
$ws4.Range("F12").Value = 296
$ws4.Range("F13").Value = 549
$ws4.Range("F14").Value = 8733
$ws4.Range("F15").Value = 7763
$ws4.Range("F16").Value = 210
$ws4.Range("F19").Value = 133
$ws4.Range("F20").Value = 992
$ws4.Range("F21").Value = 10
$ws4.Range("F26").Value = 188
$ws4.Range("F27").Value = 336
